# Fruta / hortaliza, semanal
#
# A new weekly price record for "Brocoli" at Terminal Hortofruticola Agro
# Chillan needs to be inserted as row 148 of the data table. Inserting the
# row shifts all the following rows (old 148-186) down by one (to 149-187),
# which matches the rest of the diff exactly (every downstream row keeps its
# data, just moved one row down). We then populate the brand-new row 148
# with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 148, shifting rows 148:186 down to 149:187.
$ws.Rows("148:148").Insert(-4121)  # -4121 = xlShiftDown

# Populate the newly inserted row 148 with the new data record.
$row = 148

$ws.Cells.Item($row, 1).Value2 = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"
$ws.Cells.Item($row, 4).Value2 = 44508
$ws.Cells.Item($row, 5).Value2 = 16
$ws.Cells.Item($row, 6).Value2 = 100112023
$ws.Cells.Item($row, 7).Value = "Brócoli"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value2 = 300
$ws.Cells.Item($row, 11).Value2 = 650
$ws.Cells.Item($row, 12).Value2 = 750
$ws.Cells.Item($row, 13).Value2 = 675
$ws.Cells.Item($row, 14).Value = "`$/unidad"
$ws.Cells.Item($row, 15).Value = "Región del Maule"
$ws.Cells.Item($row, 16).Value2 = 675
$ws.Cells.Item($row, 17).Value2 = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
